# Apply the "5-year growth is frozen before first drop in projections" update
# to the GroupGrowthRateMax sheet:
#   - For every 6-row block (periods 2025,2030,2035,2040,2045,2050 for one group),
#     the rate for 2035/2040/2045/2050 is frozen at the 2030 value.
#   - The notes column gets an extra clause inserted right after
#     "...logistic diffusion model" and before the "; Norway EVs growth scenario" part.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GroupGrowthRateMax")

$rowCount = $ws.UsedRange.Rows.Count

$oldPhrase = "logistic diffusion model;"
$newPhrase = "logistic diffusion model, 5-year growth is frozen before first drop in projections;"

for ($g = 2; $g -le $rowCount; $g += 6) {
    # Value to freeze forward: the 2030 rate (second row of the 6-row block)
    $frozenRate = $ws.Cells.Item($g + 1, 4).Value()

    for ($offset = 0; $offset -le 5; $offset++) {
        $r = $g + $offset

        # Freeze the rate for periods 2035, 2040, 2045, 2050 (offsets 2..5)
        if ($offset -ge 2) {
            $ws.Cells.Item($r, 4).Value = $frozenRate
        }

        # Update the notes text on every row of the block
        $notes = $ws.Cells.Item($r, 5).Value()
        if ($notes -and $notes.Contains($oldPhrase)) {
            $ws.Cells.Item($r, 5).Value = $notes.Replace($oldPhrase, $newPhrase)
        }
    }
}
